# TC07_C3DC_phs000471_DiagnBasis-NotReported.xlsx
# - Update the SurvivalTab TabQuery (cell B7) so the ORDER BY clause sorts
#   by srv.survival_id instead of prt.participant_id.
# - Re-touch the cell's font so it keeps its existing look (Calibri 12,
#   wrap text, automatic/theme text colour) but is backed by its own style.
# - Scroll the sheet back so column A is visible (topLeftCell A6 instead of B6)
#   while keeping the current selection on C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B7")

$oldText = $cell.Value2
$needle  = "    prt.participant_id ASC`nLIMIT"
$replacement = "    srv.survival_id ASC`nLIMIT"
if ($oldText.Contains($needle)) {
    $newText = $oldText.Replace($needle, $replacement)
    $cell.Value = $newText
}

# Reassert the font's theme colour; this mirrors the original file's own
# font/style (Calibri, 12pt, wrapped, automatic/theme text colour) so the
# cell keeps looking the same but is now backed by its own font/style entry.
$cell.Font.ThemeColor = 1

# Scroll the window so column A is back in view (was showing column B as the
# left-most column); the active selection (C7) is left untouched.
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
